# KIBON-120: Kanton statistic translated
# Replace hard-coded German labels with i18n template placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Parameter block ---
$ws.Range("A3").Value = "{parameterTitle}"
$ws.Range("A3").Font.Bold = $true

$ws.Range("A4").Value = "{vonTitle}"
$ws.Range("A5").Value = "{bisTitle}"

# --- Table header row 7/8 (left-to-right, then the two right-hand columns) ---
$ws.Range("A7").Value = "{fallIdTitle}"
$ws.Range("B7").Value = "{nachnameTitle}"
$ws.Range("C7").Value = "{vornameTitle}"
$ws.Range("D7").Value = "{geburtsdatumTitle}"
$ws.Range("E7").Value = "{betreuungVonTitle}"
$ws.Range("F7").Value = "{betreuungBisTitle}"
$ws.Range("G7").Value = "{bgPensumTitle}"
$ws.Range("Q7").Value = "{babyFaktorTitle}"
$ws.Range("R7").Value = "{institutionTitle}"

# --- Title row ---
$ws.Range("A1").Value = "{kantonTitle}"

# --- Remaining header row 7 cells ---
$ws.Range("H7").Value = "{monatsanfangTitle}"
$ws.Range("I7").Value = "{monatsendeTitle}"
$ws.Range("M7").Value = "{platzbelegungTageTitle}"
$ws.Range("N7").Value = "{kostenCHFTitle}"

# --- Header row 8 ---
$ws.Range("N8").Value = "{vollkostenTitle}"
$ws.Range("O8").Value = "{elternbeitragTitle}"
$ws.Range("P8").Value = "{gutscheinTitle}"

# --- Total row ---
$ws.Range("A10").Value = "{totalTitle}"

# --- Formula tweak: blank / "X" instead of Nein / Ja ---
$ws.Range("Q9").Formula = "=IF(E9>EOMONTH(D9,12),"""",""X"")"

# --- Row 8 is now taller (wrapped two-line header) ---
$ws.Rows("8").RowHeight = 30

# --- Selection moved to the header cell ---
$ws.Range("A7:A8").Select()
